# Apply the "Informe sumo" conclusions-section edit:
#  1. Extend the "La impresión funciona..." bullet with an extra clause.
#  2. Add a new bullet about printing the circuit (scaling issue).
#  3. Add a new bullet about PCB ironing/acid etching, taking over the
#     trailing _GoBack bookmark from the old empty bullet paragraph.
#  4. Refresh the cached PAGE field result in the header (1 -> 2) now
#     that the extra text pushes the document onto a second page.

$d = $word.ActiveDocument

# --- 1. "La impresión funciona..." bullet: drop the trailing period and
#        append the new clause + closing period ------------------------
$rng = $d.Content
$rng.Find.Execute(
    "soportes a la impresión.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "soportes a la impresión, garantizando así la creación efectiva de los puentes.",
    2) | Out-Null

# Locate the "La impresión funciona..." paragraph again so we can anchor
# the new paragraphs right after it (it keeps its own formatting this
# way, and the newly typed text correctly inherits the es-ES run
# formatting already present at the end of that paragraph).
$paraCount = $d.Paragraphs.Count
$idxImpresion = -1
for ($i = 1; $i -le $paraCount; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "La impresión funciona*") {
        $idxImpresion = $i
        break
    }
}
if ($idxImpresion -eq -1) {
    throw "Could not locate the 'La impresión funciona...' paragraph"
}

# --- 2. New bullet: circuit printing / scaling caution -----------------
$pImpresion = $d.Paragraphs.Item($idxImpresion)
$anchor = $d.Range($pImpresion.Range.End - 1, $pImpresion.Range.End - 1)
$anchor.InsertAfter(
    "`rPara la impresión del circuito hay que tener especial cuidado a la hora de imprimir ya que por defecto la impresión tiende a escalar dentro de la hoja, lo que provoca fallos a la hora de realizar el PCB."
) | Out-Null

# --- 3. New bullet: PCB ironing/acid-etching heat caution --------------
# This paragraph ends up owning the _GoBack bookmark that used to sit
# alone in the old empty bullet right before "Anexos".
$pCircuito = $d.Paragraphs.Item($idxImpresion + 1)
$anchor2 = $d.Range($pCircuito.Range.End - 1, $pCircuito.Range.End - 1)
$anchor2.InsertAfter(
    "`rPara la imp0’resion del PCB por el método de planchado y ácidos hay que ser cuidadoso con la cantidad de calor a impartir, para "
) | Out-Null

# --- 4. Header page-number field now reads "2" on the (now) 2nd page ---
$hdrRange = $d.Sections.Item(1).Headers.Item(1).Range
$hdrRange.Find.Execute(
    "1", $true, $false, $false, $false, $false,
    $true, 1, $false, "2", 2) | Out-Null
